$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.852.24"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "2.430.69"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.508"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("E9").Value = "  +7.15%  "
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.328"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "67.770.67"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "333.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "0.0₃0805"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "415.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.32%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0914"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("E48").Value = "  -6.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0428"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
